{"js": "// Update the \"Descriptive statistics of Republican party %votes by year\"\n// table: refreshed model run changed n, mean, sd, trimmed, skew, kurtosis,\n// se (all three data rows) plus min/range for the 2012 row.\n//\n// Table layout (row 0 is the header):\n//   0: year, 1: vars, 2: n, 3: mean, 4: sd, 5: median, 6: trimmed, 7: mad,\n//   8: min, 9: max, 10: range, 11: skew, 12: kurtosis, 13: se, 14: IQR,\n//   15: Q0.25, 16: Q0.75\n//\n// Rows: 1 -> 2008, 2 -> 2012, 3 -> 2016\n\nconst table = context.document.body.tables.getFirst();\n\nconst updates = [\n  // [rowIndex, colIndex, oldText, newText]\n  [1, 2, \"3,111\", \"3,112\"],\n  [1, 3, \"0.5695929\", \"0.5696027\"],\n  [1, 4, \"0.1630946\", \"0.1630693\"],\n  [1, 6, \"0.5726407\", \"0.5726517\"],\n  [1, 11, \"7.9695591\", \"7.9705298\"],\n  [1, 12, \"239.81650443\", \"239.8872661\"],\n  [1, 13, \"0.002924082\", \"0.002923158\"],\n\n  [2, 2, \"3,111\", \"3,112\"],\n  [2, 3, \"0.5980017\", \"0.5978287\"],\n  [2, 4, \"0.1469116\", \"0.1472042\"],\n  [2, 6, \"0.6045926\", \"0.6045143\"],\n  [2, 8, \"0.08\", \"0.06\"],\n  [2, 10, \"0.88\", \"0.90\"],\n  [2, 11, \"-0.4535041\", \"-0.4628218\"],\n  [2, 12, \"0.07026615\", \"0.1000946\"],\n  [2, 13, \"0.002633941\", \"0.002638764\"],\n\n  [3, 2, \"3,111\", \"3,112\"],\n  [3, 3, \"0.6364726\", \"0.6362936\"],\n  [3, 4, \"0.1557782\", \"0.1560723\"],\n  [3, 6, \"0.6506747\", \"0.6505778\"],\n  [3, 11, \"-0.8224014\", \"-0.8286244\"],\n  [3, 12, \"0.33974857\", \"0.3617064\"],\n  [3, 13, \"0.002792909\", \"0.002797732\"],\n];\n\nconst cells = updates.map(([r, c]) => table.getCell(r, c));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < updates.length; i++) {\n  const [, , oldText, newText] = updates[i];\n  const cell = cells[i];\n  if (cell.value !== oldText) {\n    throw new Error(\n      `Unexpected existing value \"${cell.value}\" (expected \"${oldText}\") at update index ${i}`\n    );\n  }\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the \"Descriptive statistics of Republican party %votes by year\"\n# table: refreshed model run changed n, mean, sd, trimmed, skew, kurtosis,\n# se (all three data rows) plus min/range for the 2012 row.\n#\n# Table layout (columns are 1-based for the COM object model):\n#   1: year, 2: vars, 3: n, 4: mean, 5: sd, 6: median, 7: trimmed, 8: mad,\n#   9: min, 10: max, 11: range, 12: skew, 13: kurtosis, 14: se, 15: IQR,\n#   16: Q0.25, 17: Q0.75\n#\n# Rows: 2 -> 2008, 3 -> 2012, 4 -> 2016 (row 1 is the header)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$updates = @(\n    @(2, 3, \"3,111\", \"3,112\"),\n    @(2, 4, \"0.5695929\", \"0.5696027\"),\n    @(2, 5, \"0.1630946\", \"0.1630693\"),\n    @(2, 7, \"0.5726407\", \"0.5726517\"),\n    @(2, 12, \"7.9695591\", \"7.9705298\"),\n    @(2, 13, \"239.81650443\", \"239.8872661\"),\n    @(2, 14, \"0.002924082\", \"0.002923158\"),\n\n    @(3, 3, \"3,111\", \"3,112\"),\n    @(3, 4, \"0.5980017\", \"0.5978287\"),\n    @(3, 5, \"0.1469116\", \"0.1472042\"),\n    @(3, 7, \"0.6045926\", \"0.6045143\"),\n    @(3, 9, \"0.08\", \"0.06\"),\n    @(3, 11, \"0.88\", \"0.90\"),\n    @(3, 12, \"-0.4535041\", \"-0.4628218\"),\n    @(3, 13, \"0.07026615\", \"0.1000946\"),\n    @(3, 14, \"0.002633941\", \"0.002638764\"),\n\n    @(4, 3, \"3,111\", \"3,112\"),\n    @(4, 4, \"0.6364726\", \"0.6362936\"),\n    @(4, 5, \"0.1557782\", \"0.1560723\"),\n    @(4, 7, \"0.6506747\", \"0.6505778\"),\n    @(4, 12, \"-0.8224014\", \"-0.8286244\"),\n    @(4, 13, \"0.33974857\", \"0.3617064\"),\n    @(4, 14, \"0.002792909\", \"0.002797732\")\n)\n\nforeach ($u in $updates) {\n    $rowIdx = $u[0]\n    $colIdx = $u[1]\n    $oldText = $u[2]\n    $newText = $u[3]\n    $cell = $table.Cell($rowIdx, $colIdx)\n    $range = $cell.Range\n    $range.MoveEnd(1, -1) | Out-Null\n    if ($range.Text -ne $oldText) {\n        throw \"Unexpected existing value '$($range.Text)' (expected '$oldText') at row $rowIdx col $colIdx\"\n    }\n    $range.Text = $newText\n}\n"}
